$d = $word.ActiveDocument

function Set-CellName([int]$paraIndex, [string]$name) {
    $r = $d.Paragraphs($paraIndex).Range
    $escaped = $name.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="0019302B" w:rsidRDefault="007107CA"><w:pPr><w:pStyle w:val="normal0"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="160"/></w:pPr><w:r><w:t>' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# Locate the paragraphs that currently read "Tousif" (there are four, one per
# sprint row of the Team Members column) and replace each with its correct
# team-member name, dropping the spell-check proofErr wrapper at the same time.
$names = @("SAIRAM.K", "RAJAN.V", "RAGUL.VR", "SHREE VARSHAN.R")
$targets = New-Object System.Collections.ArrayList
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Tousif") {
        [void]$targets.Add($i)
    }
}

for ($j = 0; $j -lt $targets.Count; $j++) {
    Set-CellName $targets[$j] $names[$j]
}
